$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos (USA) updated totals
$ws.Range("B4").Value = 977572
$ws.Range("C4").Value = 16921
$ws.Range("D4").Value = 118693
$ws.Range("E4").Value = 803870
$ws.Range("G4").Value = 753
$ws.Range("H4").Value = 55009

# India overtakes Peru -> rows 19/20 swap country + data
$ws.Range("A19").Value = "India"
$ws.Range("B19").Value = 27890
$ws.Range("C19").Value = 1607
$ws.Range("D19").Value = 6523
$ws.Range("E19").Value = 20485
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 57
$ws.Range("H19").Value = 882

$ws.Range("A20").Value = "Peru"
$ws.Range("B20").Value = 27517
$ws.Range("C20").Value = 2186
$ws.Range("D20").Value = 8088
$ws.Range("E20").Value = 18701
$ws.Range("F20").Value = 554
$ws.Range("G20").Value = 28
$ws.Range("H20").Value = 728

# Costa Rica (row 98)
$ws.Range("B98").Value = 695
$ws.Range("C98").Value = 2
$ws.Range("D98").Value = 264
$ws.Range("E98").Value = 425
$ws.Range("F98").Value = 8

# Sri Lanka (row 106)
$ws.Range("B106").Value = 523
$ws.Range("C106").Value = 71
$ws.Range("E106").Value = 396

# Ruanda overtakes Islas Feroe -> rows 132/133 swap country + data
$ws.Range("A132").Value = "Ruanda"
$ws.Range("B132").Value = 191
$ws.Range("C132").Value = 8
$ws.Range("D132").Value = 92
$ws.Range("E132").Value = 99

$ws.Range("A133").Value = "Islas Feroe"
$ws.Range("B133").Value = 187
$ws.Range("D133").Value = 178
$ws.Range("E133").Value = 9

# Aruba (row 148)
$ws.Range("D148").Value = 73
$ws.Range("E148").Value = 25
